$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F ("dSF") values per repulled data
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = -2
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = 5
